{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Change 1: \"Rub\u00e9n Gonzales, Diego Guzm\u00e1n, \" -> \"Rub\u00e9n Gonz\u00e1lez, Diego Guzm\u00e1n, \"\n//           (the name \"Gonzales\" is corrected to \"Gonz\u00e1lez\")\n// Change 2: \"...la venta de ropa como zapatos, pantalones, camisetas, camisas,\n//           entre otros art\u00edculos, desea ofertar...\" ->\n//           \"...la venta de art\u00edculos como: ropa, zapatos, pantalones, camisetas,\n//           camisas, entre otros. Desea ofertar...\"\n// Change 3: Remove the stray \"_GoBack\" bookmark that wrapped part of the\n//           \"Adicionalmente el vendedor...\" sentence (no visible text change,\n//           just a small syntax/cleanup correction per the commit message).\n\nconst body = context.document.body;\n\n// --- Change 1: Gonzales -> Gonz\u00e1lez -------------------------------------\nconst nameResults = body.search(\"Rub\u00e9n Gonzales\", { matchCase: true });\nnameResults.load(\"text\");\nawait context.sync();\n\nif (nameResults.items.length > 0) {\n  nameResults.items[0].insertText(\"Rub\u00e9n Gonz\u00e1lez\", \"Replace\");\n  await context.sync();\n}\n\n// --- Change 2: rewrite the \"venta de ropa como...\" sentence fragment ----\nconst oldFragment =\n  \"la venta de ropa como zapatos, pantalones, camisetas, camisas, entre otros art\u00edculos, desea ofertar\";\nconst newFragment =\n  \"la venta de art\u00edculos como: ropa, zapatos, pantalones, camisetas, camisas, entre otros. Desea ofertar\";\n\nconst fragmentResults = body.search(oldFragment, { matchCase: true });\nfragmentResults.load(\"text\");\nawait context.sync();\n\nif (fragmentResults.items.length > 0) {\n  fragmentResults.items[0].insertText(newFragment, \"Replace\");\n  await context.sync();\n}\n\n// --- Change 3: remove the \"_GoBack\" bookmark and reunify the sentence ---\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\nconst sentence =\n  \"Adicionalmente el vendedor desea obtener la estad\u00edstica para conocer el o los productos m\u00e1s vendidos.\";\nconst sentenceResults = body.search(sentence, { matchCase: true });\nsentenceResults.load(\"text\");\nawait context.sync();\n\nif (sentenceResults.items.length > 0) {\n  sentenceResults.items[0].insertText(sentence, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word / $d ($word.ActiveDocument) are pre-seeded by the host.\n#\n# Change 1: \"Rub\u00e9n Gonzales, Diego Guzm\u00e1n, \" -> \"Rub\u00e9n Gonz\u00e1lez, Diego Guzm\u00e1n, \"\n#           (the name \"Gonzales\" is corrected to \"Gonz\u00e1lez\")\n# Change 2: \"...la venta de ropa como zapatos, pantalones, camisetas, camisas,\n#           entre otros art\u00edculos, desea ofertar...\" ->\n#           \"...la venta de art\u00edculos como: ropa, zapatos, pantalones, camisetas,\n#           camisas, entre otros. Desea ofertar...\"\n# Change 3: Remove the stray \"_GoBack\" bookmark that wrapped part of the\n#           \"Adicionalmente el vendedor...\" sentence (no visible text change,\n#           just a small syntax/cleanup correction per the commit message).\n\n$d = $word.ActiveDocument\n\nfunction Replace-Text($findText, $replaceText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    # wdReplaceOne (1) replaces only the first match found.\n    $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $false, $false, $find.Forward, $find.Wrap, $false, $find.Replacement.Text, 1) | Out-Null\n}\n\n# --- Change 1: Gonzales -> Gonz\u00e1lez -------------------------------------\nReplace-Text \"Rub\u00e9n Gonzales\" \"Rub\u00e9n Gonz\u00e1lez\"\n\n# --- Change 2: rewrite the \"venta de ropa como...\" sentence fragment ----\n$oldFragment = \"la venta de ropa como zapatos, pantalones, camisetas, camisas, entre otros art\u00edculos, desea ofertar\"\n$newFragment = \"la venta de art\u00edculos como: ropa, zapatos, pantalones, camisetas, camisas, entre otros. Desea ofertar\"\nReplace-Text $oldFragment $newFragment\n\n# --- Change 3: remove the \"_GoBack\" bookmark and reunify the sentence ---\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n\n$sentence = \"Adicionalmente el vendedor desea obtener la estad\u00edstica para conocer el o los productos m\u00e1s vendidos.\"\nReplace-Text $sentence $sentence\n"}
